# Sprint 4 burndown chart — update the "Actual" daily burn row (row 23)
# with the real numbers for the sprint; the dependent running-total row
# (24) and the chart series simply recalc off of these via their
# existing formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 4")
$ws.Activate()

$actual = @{
    "D23" = 2
    "E23" = 7
    "F23" = 5
    "G23" = 7
    "H23" = 4
    "I23" = 5
    "J23" = 6
    "K23" = 3
    "L23" = 2
    "M23" = 3
    "N23" = 3
    "O23" = 2
    "P23" = 3
    "Q23" = 2
}

foreach ($addr in $actual.Keys) {
    $ws.Range($addr).Value = $actual[$addr]
}

# Match the saved selection/scroll state from the edit: the author had
# scrolled down a few rows and had H23 selected when they saved.
$excel.Goto($ws.Range("A14"), $true)
$ws.Range("H23").Select()
